$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 921.63635
$ws.Range("I19").Value = 809.8570999999999
$ws.Range("K19").Value = 809.8570999999999
$ws.Range("M19").Value = -634.8570999999999
$ws.Range("H40").Value = 2057.2
$ws.Range("I40").Value = 1532.6666
$ws.Range("J40").Value = 2844
$ws.Range("K40").Value = 1532.6666
$ws.Range("L40").Value = 2844
$ws.Range("M40").Value = -1357.6666
$ws.Range("N40").Value = -3194
$ws.Range("H62").Value = 4999.5
$ws.Range("I62").Value = 4999.5
$ws.Range("K62").Value = 4999.5
$ws.Range("M62").Value = -4375.5
$ws.Range("H65").Value = 4999.5
$ws.Range("I65").Value = 4999.5
$ws.Range("K65").Value = 24997.5
$ws.Range("M65").Value = -21877.5
$ws.Range("H103").Value = 345.66666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 594.375
$ws.Range("J2").Value = 865.6667
$ws.Range("L2").Value = 865.6667
$ws.Range("N2").Value = -1091.6667
$ws.Range("H61").Value = 3208991.8
$ws.Range("I61").Value = 67489.234
$ws.Range("K61").Value = 67489.234
$ws.Range("M61").Value = -67277.234
$ws.Range("H63").Value = 21757.812
$ws.Range("I63").Value = 7708.6665
$ws.Range("K63").Value = 7708.6665
$ws.Range("M63").Value = -7022.6665
$ws.Range("H66").Value = 21757.812
$ws.Range("I66").Value = 7708.6665
$ws.Range("K66").Value = 38543.3325
$ws.Range("M66").Value = -35111.3325
$ws.Range("H74").Value = 708536.0600000001
$ws.Range("I74").Value = 4944.174
$ws.Range("K74").Value = 4944.174
$ws.Range("M74").Value = -4070.174
$ws.Range("H77").Value = 708536.0600000001
$ws.Range("I77").Value = 4944.174
$ws.Range("K77").Value = 24720.87
$ws.Range("M77").Value = -20352.87
$ws.Range("H116").Value = 594.375
$ws.Range("J116").Value = 865.6667
$ws.Range("L116").Value = 865.6667
$ws.Range("N116").Value = -5453.6667
$ws.Range("H132").Value = 900
$ws.Range("I132").Value = 900
$ws.Range("K132").Value = 2700
$ws.Range("M132").Value = -170
$ws.Range("H136").Value = 3208991.8
$ws.Range("I136").Value = 67489.234
$ws.Range("K136").Value = 202467.702
$ws.Range("M136").Value = -199917.702

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 594.375
$ws.Range("J3").Value = 865.6667
$ws.Range("L3").Value = 865.6667
$ws.Range("N3").Value = -1093.6667
$ws.Range("H7").Value = 2525000
$ws.Range("I7").Value = 2525000
$ws.Range("K7").Value = 2525000
$ws.Range("M7").Value = -2524887

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2799.0977
$ws.Range("I31").Value = 3411.5264
$ws.Range("K31").Value = 3411.5264
$ws.Range("M31").Value = -3116.5264
$ws.Range("H34").Value = 2799.0977
$ws.Range("I34").Value = 3411.5264
$ws.Range("K34").Value = 3411.5264
$ws.Range("M34").Value = -3209.5264
$ws.Range("H94").Value = 1195
$ws.Range("J94").Value = 1175.125
$ws.Range("L94").Value = 1175.125
$ws.Range("N94").Value = -2077.125
$ws.Range("H100").Value = 20899.5
$ws.Range("J100").Value = 20899.5
$ws.Range("L100").Value = 20899.5
$ws.Range("N100").Value = -23063.5
$ws.Range("H132").Value = 55996.844
$ws.Range("I132").Value = 94302.09
$ws.Range("J132").Value = 3327.125
$ws.Range("K132").Value = 282906.27
$ws.Range("L132").Value = 9981.375
$ws.Range("M132").Value = -280376.27
$ws.Range("N132").Value = -15041.375
$ws.Range("H134").Value = 2221.647
$ws.Range("I134").Value = 2083.44
$ws.Range("J134").Value = 2605.5557
$ws.Range("K134").Value = 6250.32
$ws.Range("L134").Value = 7816.6671
$ws.Range("M134").Value = -3715.32
$ws.Range("N134").Value = -12886.6671

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 236.4
$ws.Range("I14").Value = 236.4
$ws.Range("K14").Value = 709.2
$ws.Range("M14").Value = -536.2
$ws.Range("H64").Value = 8711.857
$ws.Range("I64").Value = 5329.6665
$ws.Range("J64").Value = 11248.5
$ws.Range("K64").Value = 15988.9995
$ws.Range("L64").Value = 33745.5
$ws.Range("M64").Value = -15718.9995
$ws.Range("N64").Value = -34285.5
$ws.Range("H67").Value = 8711.857
$ws.Range("I67").Value = 5329.6665
$ws.Range("J67").Value = 11248.5
$ws.Range("K67").Value = 15988.9995
$ws.Range("L67").Value = 33745.5
$ws.Range("M67").Value = -15052.9995
$ws.Range("N67").Value = -35617.5
$ws.Range("H75").Value = 27783402
$ws.Range("J75").Value = 35721428
$ws.Range("L75").Value = 107164284
$ws.Range("N75").Value = -107166280
$ws.Range("H78").Value = 27783402
$ws.Range("J78").Value = 35721428
$ws.Range("L78").Value = 321492852
$ws.Range("N78").Value = -321502836
$ws.Range("H104").Value = 10175.8
$ws.Range("I104").Value = 4000
$ws.Range("J104").Value = 10862
$ws.Range("K104").Value = 12000
$ws.Range("L104").Value = 32586
$ws.Range("M104").Value = -9379
$ws.Range("N104").Value = -37828
$ws.Range("H107").Value = 1090.5883
$ws.Range("I107").Value = 1149.5
$ws.Range("K107").Value = 3448.5
$ws.Range("M107").Value = -1528.5
$ws.Range("H113").Value = 515.44446
$ws.Range("I113").Value = 617
$ws.Range("J113").Value = 434.2
$ws.Range("K113").Value = 1851
$ws.Range("L113").Value = 1302.6
$ws.Range("M113").Value = 319
$ws.Range("N113").Value = -5642.6
$ws.Range("H133").Value = 1999
$ws.Range("I133").Value = 1999
$ws.Range("K133").Value = 5997
$ws.Range("M133").Value = -937
$ws.Range("H141").Value = 1516.2
$ws.Range("I141").Value = 1516.2
$ws.Range("K141").Value = 4548.6
$ws.Range("M141").Value = 631.3999999999996

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 26500
$ws.Range("H113").Value = 2377.75
$ws.Range("I113").Value = 1011
$ws.Range("K113").Value = 1011
$ws.Range("M113").Value = 1159
$ws.Range("H136").Value = 37998.25
$ws.Range("J136").Value = 37998.25
$ws.Range("L136").Value = 113994.75
$ws.Range("N136").Value = -119094.75
$ws.Range("H139").Value = 156098.5
$ws.Range("J139").Value = 176498.75
$ws.Range("L139").Value = 176498.75
$ws.Range("N139").Value = -186778.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7735.5713
$ws.Range("I7").Value = 3480.4443
$ws.Range("K7").Value = 3480.4443
$ws.Range("M7").Value = -3368.4443
$ws.Range("H40").Value = 5956.7144
$ws.Range("I40").Value = 5939.4
$ws.Range("K40").Value = 5939.4
$ws.Range("M40").Value = -5803.4
$ws.Range("H43").Value = 1021525.75
$ws.Range("J43").Value = 1261999.4
$ws.Range("L43").Value = 1261999.4
$ws.Range("N43").Value = -1262385.4
$ws.Range("H61").Value = 2112.8125
$ws.Range("I61").Value = 2081.7856
$ws.Range("K61").Value = 2081.7856
$ws.Range("M61").Value = -1879.7856
$ws.Range("H113").Value = 2112.8125
$ws.Range("I113").Value = 2081.7856
$ws.Range("K113").Value = 2081.7856
$ws.Range("M113").Value = 88.21439999999984
$ws.Range("H126").Value = 7735.5713
$ws.Range("I126").Value = 3480.4443
$ws.Range("K126").Value = 10441.3329
$ws.Range("M126").Value = -7971.332900000001
$ws.Range("H127").Value = 95281.86
$ws.Range("J127").Value = 95281.86
$ws.Range("L127").Value = 95281.86
$ws.Range("N127").Value = -105201.86
$ws.Range("H132").Value = 3297.963
$ws.Range("I132").Value = 3193.6667
$ws.Range("K132").Value = 9581.000100000001
$ws.Range("M132").Value = -7051.000100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 35000
$ws.Range("I64").Value = 35000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 35000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -34752
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 35000
$ws.Range("I67").Value = 35000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 35000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -34142
$ws.Range("N67").ClearContents()
$ws.Range("H100").Value = 834
$ws.Range("I100").Value = 528.8333
$ws.Range("K100").Value = 1057.6666
$ws.Range("M100").Value = -516.6666
$ws.Range("H113").Value = 363.8889
$ws.Range("I113").Value = 107.333336
$ws.Range("J113").Value = 492.16666
$ws.Range("K113").Value = 322.000008
$ws.Range("L113").Value = 1476.49998
$ws.Range("M113").Value = 1847.999992
$ws.Range("N113").Value = -5816.499980000001
$ws.Range("H132").Value = 2676.2083
$ws.Range("I132").Value = 2389.842
$ws.Range("J132").Value = 3764.4
$ws.Range("K132").Value = 7169.526
$ws.Range("L132").Value = 11293.2
$ws.Range("M132").Value = -4639.526
$ws.Range("N132").Value = -16353.2
